$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new day's figure (date=2020-04-12 serial 43933, deaths=103) was added
# to the top of the data table. Insert a fresh row right under the header
# (row 1) so every existing data row shifts down by one.
# Use CopyOrigin = xlFormatFromRightOrBelow (-4121) so the new row doesn't
# inherit the bold/centered header formatting from row 1 above it.
$ws.Rows.Item(2).Insert(-4121, 0)

# Reset formatting picked up from the insert back to the sheet's default,
# then re-apply the date number format used by the rest of column A
# (the column style already carries it, but set it explicitly to be safe).
$ws.Cells.Item(2, 1).Style = "Normal"
$ws.Cells.Item(2, 2).Style = "Normal"
$ws.Cells.Item(2, 1).NumberFormat = "yyyy-mm-dd"

# Populate the new row with the latest date/deaths entry
$ws.Cells.Item(2, 1).Value = 43933
$ws.Cells.Item(2, 2).Value = 103
